# Update the table style used by the table on the "Cash flow" slide
# (Slide 16) from the deck's default table style to the built-in
# PowerPoint table style {8814871F-9231-4903-8F3B-C9315190C464}.

$p = $ppt.ActivePresentation

$targetStyleId = "{8814871F-9231-4903-8F3B-C9315190C464}"

foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($targetStyleId)
        }
    }
}
